$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting date/condition columns
# (B:E) one column to the right (C:F). This mirrors the recorded change
# of moving data recording earlier (new "BA" / numeric-date columns).
$ws.Columns("B:B").Insert()

# Record the new data in the freshly inserted column B and existing column C
$ws.Range("B2").Value = "BA"
$ws.Range("C2").Value = 20250107

# Update the active selection to reflect where the user left off (G2)
$ws.Range("G2").Select()
